$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.108.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.827.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +7.86%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.86%  "
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.825.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.81%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.342.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.016.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.813.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  +8.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000104"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "512.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "184.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.13%  "
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0856"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.573"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.91%  "
$ws.Range("E51").Value = "  +3.00%  "
